# Replace annual S&L tax revenue with quarterly numbers
# Updates columns E, F, G, H, I for rows 87-93 (2021 Q1 - 2022 Q3)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E87").Value = 3.274442939439371
$ws.Range("F87").Value = 4.485867473255149
$ws.Range("G87").Value = 3.100508840408109
$ws.Range("H87").Value = -3.368956404567739
$ws.Range("I87").Value = 3.542890503599001

$ws.Range("E88").Value = -2.987287850227718
$ws.Range("F88").Value = 0.1666011145068151
$ws.Range("G88").Value = -0.7037459759121391
$ws.Range("H88").Value = 0.2868626796071934
$ws.Range("I88").Value = -2.570404553922772

$ws.Range("E89").Value = -2.924166013534065
$ws.Range("F89").Value = -1.486200667880462
$ws.Range("G89").Value = -0.1499749553841785
$ws.Range("H89").Value = -0.1256411816786448
$ws.Range("I89").Value = -2.648549876471242

$ws.Range("E90").Value = -2.597578075119432
$ws.Range("F90").Value = -1.308647249860461
$ws.Range("G90").Value = -3.554553272110574
$ws.Range("H90").Value = 3.433091938529865
$ws.Range("I90").Value = -2.476116741538723

$ws.Range("E91").Value = -6.368117040881671
$ws.Range("F91").Value = -3.719287244940721
$ws.Range("I91").Value = -6.311287036883249

$ws.Range("E92").Value = -1.976178782683645
$ws.Range("F92").Value = -3.466509978054704
$ws.Range("I92").Value = -1.918836563185962

$ws.Range("E93").Value = -0.5269966399726116
$ws.Range("F93").Value = -2.867217634664341
$ws.Range("I93").Value = -0.4317075462883864
